$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting so that values such as
# "1.000" or "26.798.04" are not re-interpreted as numbers/dates by Excel.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '26.798.04'
$ws.Range('E2').Value = '  -3.10%  '
$ws.Range('D3').Value = '1.855.03'
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '304.69'
$ws.Range('E5').Value = '  -1.75%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '0.5087'
$ws.Range('E7').Value = '  -3.47%  '
$ws.Range('D8').Value = '0.3653'
$ws.Range('E8').Value = '  -4.03%  '
$ws.Range('D9').Value = '0.07136'
$ws.Range('E9').Value = '  -1.37%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = '20.72'
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = '0.8867'
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('D12').Value = '0.07513'
$ws.Range('E12').Value = '  -1.58%  '
$ws.Range('D13').Value = '1.859.51'
$ws.Range('E13').Value = '  -1.81%  '
$ws.Range('D14').Value = '5.239'
$ws.Range('E14').Value = '  -3.52%  '
$ws.Range('D15').Value = '91.03'
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '0.000008530'
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('D18').Value = '14.05'
$ws.Range('E18').Value = '  -1.89%  '
$ws.Range('D19').Value = '0.9999'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = '26.854.38'
$ws.Range('E20').Value = '  -3.02%  '
$ws.Range('D21').Value = '5.005'
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('D22').Value = '2.089.96'
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('D23').Value = '10.24'
$ws.Range('E23').Value = '  -5.21%  '
$ws.Range('D24').Value = '6.437'
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('D25').Value = '1.821'
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('D26').Value = '146.17'
$ws.Range('E26').Value = '  -4.71%  '
$ws.Range('D27').Value = '17.82'
$ws.Range('E27').Value = '  -2.17%  '
$ws.Range('D28').Value = '2.047'
$ws.Range('E28').Value = '  -6.81%  '
$ws.Range('D29').Value = '112.88'
$ws.Range('E29').Value = '  -1.11%  '
$ws.Range('D30').Value = '4.628'
$ws.Range('E30').Value = '  -4.27%  '
$ws.Range('D31').Value = '4.668'
$ws.Range('E31').Value = '  -2.85%  '
$ws.Range('D32').Value = '0.09259'
$ws.Range('E32').Value = '  +1.24%  '
$ws.Range('D33').Value = '0.05104'
$ws.Range('E33').Value = '  -3.14%  '
$ws.Range('D34').Value = '3.067'
$ws.Range('E34').Value = '  -1.76%  '
$ws.Range('E35').Value = '  -5.68%  '
$ws.Range('D36').Value = '0.7308'
$ws.Range('E36').Value = '  -5.14%  '
$ws.Range('D37').Value = '3.188'
$ws.Range('E37').Value = '  +3.64%  '
$ws.Range('D38').Value = '0.02008'
$ws.Range('E38').Value = '  -3.65%  '
$ws.Range('D39').Value = '2.464'
$ws.Range('E39').Value = '  -3.92%  '
$ws.Range('D40').Value = '1.074'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('D41').Value = '0.5280'
$ws.Range('E41').Value = '  -5.33%  '
$ws.Range('D42').Value = '117.39'
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('D43').Value = '6.463'
$ws.Range('E43').Value = '  -3.78%  '
$ws.Range('D44').Value = '8.371'
$ws.Range('E44').Value = '  -3.68%  '
$ws.Range('E45').Value = '  -2.45%  '
$ws.Range('D46').Value = '0.4644'
$ws.Range('E46').Value = '  -3.33%  '
$ws.Range('D47').Value = '0.9997'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('E48').Value = '  -4.98%  '
$ws.Range('D49').Value = '1.554'
$ws.Range('E49').Value = '  -2.61%  '
$ws.Range('D50').Value = '36.95'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').Value = '62.94'
$ws.Range('E51').Value = '  -4.97%  '
